$wb = $excel.ActiveWorkbook

# ---- workbook.xml: update the hidden _xlnm._FilterDatabase name for the Attendance sheet
#      so its range keeps pace with the newly appended rows (A1:K688 -> A1:K712) ----
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -eq "Attendance!_FilterDatabase") {
        $nm.RefersTo = "='Attendance'!`$A`$1:`$K`$712"
    }
}

# ---- Summary sheet: update attendance stats for students with a new MICROBIOLOGY session ----
$summary = $wb.Worksheets.Item("Summary")

# Row 23
$cell = $summary.Range("G23")
$cell.NumberFormat = "@"
$cell.Value = "13.8%"
$cell.NumberFormat = "0.0%"
$summary.Range("I23").Value = 19
$summary.Range("N23").Value = 4
$summary.Range("O23").Value = 8
$summary.Range("AG23").Value = 1

# Row 38
$cell = $summary.Range("G38")
$cell.NumberFormat = "@"
$cell.Value = "24.1%"
$cell.NumberFormat = "0.0%"
$summary.Range("I38").Value = 16
$summary.Range("N38").Value = 7
$summary.Range("O38").Value = 5
$summary.Range("AG38").Value = 1

# Row 55
$cell = $summary.Range("G55")
$cell.NumberFormat = "@"
$cell.Value = "13.8%"
$cell.NumberFormat = "0.0%"
$summary.Range("I55").Value = 19
$summary.Range("N55").Value = 4
$summary.Range("O55").Value = 8
$summary.Range("AG55").Value = 1

# Row 100
$cell = $summary.Range("G100")
$cell.NumberFormat = "@"
$cell.Value = "17.2%"
$cell.NumberFormat = "0.0%"
$summary.Range("I100").Value = 18
$summary.Range("N100").Value = 5
$summary.Range("O100").Value = 7
$summary.Range("AG100").Value = 1

# Row 108
$cell = $summary.Range("G108")
$cell.NumberFormat = "@"
$cell.Value = "17.2%"
$cell.NumberFormat = "0.0%"
$summary.Range("I108").Value = 18
$summary.Range("N108").Value = 5
$summary.Range("O108").Value = 7
$summary.Range("AG108").Value = 1

# Row 110
$cell = $summary.Range("G110")
$cell.NumberFormat = "@"
$cell.Value = "17.2%"
$cell.NumberFormat = "0.0%"
$summary.Range("I110").Value = 18
$summary.Range("N110").Value = 5
$summary.Range("O110").Value = 7
$summary.Range("AG110").Value = 1

# Row 111
$cell = $summary.Range("G111")
$cell.NumberFormat = "@"
$cell.Value = "17.2%"
$cell.NumberFormat = "0.0%"
$summary.Range("I111").Value = 18
$summary.Range("N111").Value = 5
$summary.Range("O111").Value = 7
$summary.Range("AG111").Value = 1

# Row 130
$cell = $summary.Range("G130")
$cell.NumberFormat = "@"
$cell.Value = "24.1%"
$cell.NumberFormat = "0.0%"
$summary.Range("I130").Value = 16
$summary.Range("N130").Value = 7
$summary.Range("O130").Value = 5
$summary.Range("AG130").Value = 1

# Row 134
$cell = $summary.Range("F134")
$summary.Range("F3").Copy()
$cell.PasteSpecial(-4122)
$cell.Value = "Moderate Risk"
$cell = $summary.Range("G134")
$cell.NumberFormat = "@"
$cell.Value = "27.6%"
$cell.NumberFormat = "0.0%"
$summary.Range("I134").Value = 15
$summary.Range("N134").Value = 8
$summary.Range("O134").Value = 4
$summary.Range("AG134").Value = 1

# Row 143
$cell = $summary.Range("G143")
$cell.NumberFormat = "@"
$cell.Value = "17.2%"
$cell.NumberFormat = "0.0%"
$summary.Range("I143").Value = 18
$summary.Range("N143").Value = 5
$summary.Range("O143").Value = 7
$summary.Range("AG143").Value = 1

# Row 145
$cell = $summary.Range("F145")
$summary.Range("F2").Copy()
$cell.PasteSpecial(-4122)
$cell.Value = "High Risk"
$cell = $summary.Range("G145")
$cell.NumberFormat = "@"
$cell.Value = "20.7%"
$cell.NumberFormat = "0.0%"
$summary.Range("I145").Value = 17
$summary.Range("N145").Value = 6
$summary.Range("O145").Value = 6
$summary.Range("AG145").Value = 1

# Row 146
$cell = $summary.Range("F146")
$summary.Range("F2").Copy()
$cell.PasteSpecial(-4122)
$cell.Value = "High Risk"
$cell = $summary.Range("G146")
$cell.NumberFormat = "@"
$cell.Value = "20.7%"
$cell.NumberFormat = "0.0%"
$summary.Range("I146").Value = 17
$summary.Range("N146").Value = 6
$summary.Range("O146").Value = 6
$summary.Range("AG146").Value = 1

# Row 148
$cell = $summary.Range("G148")
$cell.NumberFormat = "@"
$cell.Value = "17.2%"
$cell.NumberFormat = "0.0%"
$summary.Range("I148").Value = 18
$summary.Range("N148").Value = 5
$summary.Range("O148").Value = 7
$summary.Range("AG148").Value = 1

# Row 150
$cell = $summary.Range("G150")
$cell.NumberFormat = "@"
$cell.Value = "10.3%"
$cell.NumberFormat = "0.0%"
$summary.Range("I150").Value = 20
$summary.Range("N150").Value = 3
$summary.Range("O150").Value = 9
$summary.Range("AG150").Value = 1

# Row 151
$cell = $summary.Range("G151")
$cell.NumberFormat = "@"
$cell.Value = "10.3%"
$cell.NumberFormat = "0.0%"
$summary.Range("I151").Value = 20
$summary.Range("N151").Value = 3
$summary.Range("O151").Value = 9
$summary.Range("AG151").Value = 1

# Row 152
$cell = $summary.Range("F152")
$summary.Range("F2").Copy()
$cell.PasteSpecial(-4122)
$cell.Value = "High Risk"
$cell = $summary.Range("G152")
$cell.NumberFormat = "@"
$cell.Value = "20.7%"
$cell.NumberFormat = "0.0%"
$summary.Range("I152").Value = 17
$summary.Range("N152").Value = 6
$summary.Range("O152").Value = 6
$summary.Range("AG152").Value = 1

# Row 157
$cell = $summary.Range("G157")
$cell.NumberFormat = "@"
$cell.Value = "10.3%"
$cell.NumberFormat = "0.0%"
$summary.Range("I157").Value = 20
$summary.Range("N157").Value = 3
$summary.Range("O157").Value = 9
$summary.Range("AG157").Value = 1

# Row 175
$cell = $summary.Range("F175")
$summary.Range("F2").Copy()
$cell.PasteSpecial(-4122)
$cell.Value = "High Risk"
$cell = $summary.Range("G175")
$cell.NumberFormat = "@"
$cell.Value = "20.7%"
$cell.NumberFormat = "0.0%"
$summary.Range("I175").Value = 17
$summary.Range("N175").Value = 6
$summary.Range("O175").Value = 6
$summary.Range("AG175").Value = 1

# Row 193
$cell = $summary.Range("G193")
$cell.NumberFormat = "@"
$cell.Value = "10.3%"
$cell.NumberFormat = "0.0%"
$summary.Range("I193").Value = 20
$summary.Range("N193").Value = 3
$summary.Range("O193").Value = 9
$summary.Range("AG193").Value = 1

# Row 194
$cell = $summary.Range("G194")
$cell.NumberFormat = "@"
$cell.Value = "10.3%"
$cell.NumberFormat = "0.0%"
$summary.Range("I194").Value = 20
$summary.Range("N194").Value = 3
$summary.Range("O194").Value = 9
$summary.Range("AG194").Value = 1

# Row 198
$cell = $summary.Range("G198")
$cell.NumberFormat = "@"
$cell.Value = "13.8%"
$cell.NumberFormat = "0.0%"
$summary.Range("I198").Value = 19
$summary.Range("N198").Value = 4
$summary.Range("O198").Value = 8
$summary.Range("AG198").Value = 1

# Row 228
$cell = $summary.Range("G228")
$cell.NumberFormat = "@"
$cell.Value = "17.2%"
$cell.NumberFormat = "0.0%"
$summary.Range("I228").Value = 18
$summary.Range("N228").Value = 5
$summary.Range("O228").Value = 7
$summary.Range("AG228").Value = 1

# Row 241
$cell = $summary.Range("F241")
$summary.Range("F3").Copy()
$cell.PasteSpecial(-4122)
$cell.Value = "Moderate Risk"
$cell = $summary.Range("G241")
$cell.NumberFormat = "@"
$cell.Value = "27.6%"
$cell.NumberFormat = "0.0%"
$summary.Range("I241").Value = 15
$summary.Range("N241").Value = 8
$summary.Range("O241").Value = 4
$summary.Range("AG241").Value = 1

# Row 252
$cell = $summary.Range("F252")
$summary.Range("F2").Copy()
$cell.PasteSpecial(-4122)
$cell.Value = "High Risk"
$cell = $summary.Range("G252")
$cell.NumberFormat = "@"
$cell.Value = "20.7%"
$cell.NumberFormat = "0.0%"
$summary.Range("I252").Value = 17
$summary.Range("N252").Value = 6
$summary.Range("O252").Value = 6
$summary.Range("AG252").Value = 1

$excel.CutCopyMode = $false

# ---- Attendance sheet: append 24 new MICROBIOLOGY attendance rows (689-712) ----
$att = $wb.Worksheets.Item("Attendance")
$plainTextSrc = $att.Range("A2")   # a cell with default (unstyled) format, used to strip styling from new text cells

# Row 689 - student 221031
$c = $att.Cells.Item(689,1)
$c.NumberFormat = "@"
$c.Value = "221031"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(689,2).Value = "امنيه عبدالله عبد اللطيف محمد"
$att.Cells.Item(689,3).Value = "Year 2"
$att.Cells.Item(689,4).Value = "C1"
$att.Cells.Item(689,5).Value = "221031@med.asu.edu.eg"
$att.Cells.Item(689,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(689,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(689,8).Value = "MICROBIOLOGY"
$att.Cells.Item(689,9).Value = "24/11/2025"
$att.Cells.Item(689,10).Value = "09:21:14"
$att.Cells.Item(689,11).Value = "C1"

# Row 690 - student 221584
$c = $att.Cells.Item(690,1)
$c.NumberFormat = "@"
$c.Value = "221584"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(690,2).Value = "عزه بنت محمد بن عوض الصمداني"
$att.Cells.Item(690,3).Value = "Year 2"
$att.Cells.Item(690,4).Value = "C1"
$att.Cells.Item(690,5).Value = "221584@med.asu.edu.eg"
$att.Cells.Item(690,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(690,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(690,8).Value = "MICROBIOLOGY"
$att.Cells.Item(690,9).Value = "24/11/2025"
$att.Cells.Item(690,10).Value = "09:21:33"
$att.Cells.Item(690,11).Value = "C1"

# Row 691 - student 221307
$c = $att.Cells.Item(691,1)
$c.NumberFormat = "@"
$c.Value = "221307"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(691,2).Value = "دعاء عاصم على العوض"
$att.Cells.Item(691,3).Value = "Year 2"
$att.Cells.Item(691,4).Value = "C1"
$att.Cells.Item(691,5).Value = "221307@med.asu.edu.eg"
$att.Cells.Item(691,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(691,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(691,8).Value = "MICROBIOLOGY"
$att.Cells.Item(691,9).Value = "24/11/2025"
$att.Cells.Item(691,10).Value = "09:21:41"
$att.Cells.Item(691,11).Value = "C1"

# Row 692 - student 221818
$c = $att.Cells.Item(692,1)
$c.NumberFormat = "@"
$c.Value = "221818"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(692,2).Value = "جيهان محارب الشيخ الكيلاني"
$att.Cells.Item(692,3).Value = "Year 2"
$att.Cells.Item(692,4).Value = "C1"
$att.Cells.Item(692,5).Value = "221818@med.asu.edu.eg"
$att.Cells.Item(692,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(692,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(692,8).Value = "MICROBIOLOGY"
$att.Cells.Item(692,9).Value = "24/11/2025"
$att.Cells.Item(692,10).Value = "09:21:52"
$att.Cells.Item(692,11).Value = "C1"

# Row 693 - student 221810
$c = $att.Cells.Item(693,1)
$c.NumberFormat = "@"
$c.Value = "221810"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(693,2).Value = "رهان محارب الشيخ الكيلاني"
$att.Cells.Item(693,3).Value = "Year 2"
$att.Cells.Item(693,4).Value = "C1"
$att.Cells.Item(693,5).Value = "221810@med.asu.edu.eg"
$att.Cells.Item(693,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(693,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(693,8).Value = "MICROBIOLOGY"
$att.Cells.Item(693,9).Value = "24/11/2025"
$att.Cells.Item(693,10).Value = "09:22:09"
$att.Cells.Item(693,11).Value = "C1"

# Row 694 - student 221838
$c = $att.Cells.Item(694,1)
$c.NumberFormat = "@"
$c.Value = "221838"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(694,2).Value = "الاء سيد احمد احمد ابن ادريس"
$att.Cells.Item(694,3).Value = "Year 2"
$att.Cells.Item(694,4).Value = "C1"
$att.Cells.Item(694,5).Value = "221838@med.asu.edu.eg"
$att.Cells.Item(694,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(694,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(694,8).Value = "MICROBIOLOGY"
$att.Cells.Item(694,9).Value = "24/11/2025"
$att.Cells.Item(694,10).Value = "09:22:20"
$att.Cells.Item(694,11).Value = "C1"

# Row 695 - student 222035
$c = $att.Cells.Item(695,1)
$c.NumberFormat = "@"
$c.Value = "222035"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(695,2).Value = "علا عبد الوهاب خليل محمود"
$att.Cells.Item(695,3).Value = "Year 2"
$att.Cells.Item(695,4).Value = "C1"
$att.Cells.Item(695,5).Value = "222035@med.asu.edu.eg"
$att.Cells.Item(695,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(695,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(695,8).Value = "MICROBIOLOGY"
$att.Cells.Item(695,9).Value = "24/11/2025"
$att.Cells.Item(695,10).Value = "09:22:30"
$att.Cells.Item(695,11).Value = "C1"

# Row 696 - student 221319
$c = $att.Cells.Item(696,1)
$c.NumberFormat = "@"
$c.Value = "221319"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(696,2).Value = "روان صلاح طاهر الوهباني"
$att.Cells.Item(696,3).Value = "Year 2"
$att.Cells.Item(696,4).Value = "C1"
$att.Cells.Item(696,5).Value = "221319@med.asu.edu.eg"
$att.Cells.Item(696,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(696,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(696,8).Value = "MICROBIOLOGY"
$att.Cells.Item(696,9).Value = "24/11/2025"
$att.Cells.Item(696,10).Value = "09:22:39"
$att.Cells.Item(696,11).Value = "C1"

# Row 697 - student 210967
$c = $att.Cells.Item(697,1)
$c.NumberFormat = "@"
$c.Value = "210967"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(697,2).Value = "ملاك كمال اسماعيل ابو جلاله"
$att.Cells.Item(697,3).Value = "Year 2"
$att.Cells.Item(697,4).Value = "C1"
$att.Cells.Item(697,5).Value = "210967@med.asu.edu.eg"
$att.Cells.Item(697,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(697,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(697,8).Value = "MICROBIOLOGY"
$att.Cells.Item(697,9).Value = "24/11/2025"
$att.Cells.Item(697,10).Value = "09:23:01"
$att.Cells.Item(697,11).Value = "C1"

# Row 698 - student 222113
$c = $att.Cells.Item(698,1)
$c.NumberFormat = "@"
$c.Value = "222113"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(698,2).Value = "ابرار محمد عبد الله عبد الحميد"
$att.Cells.Item(698,3).Value = "Year 2"
$att.Cells.Item(698,4).Value = "C1"
$att.Cells.Item(698,5).Value = "222113@med.asu.edu.eg"
$att.Cells.Item(698,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(698,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(698,8).Value = "MICROBIOLOGY"
$att.Cells.Item(698,9).Value = "24/11/2025"
$att.Cells.Item(698,10).Value = "09:23:10"
$att.Cells.Item(698,11).Value = "C1"

# Row 699 - student 221675
$c = $att.Cells.Item(699,1)
$c.NumberFormat = "@"
$c.Value = "221675"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(699,2).Value = "ساره بنت سعيد بن عثمان الكناني"
$att.Cells.Item(699,3).Value = "Year 2"
$att.Cells.Item(699,4).Value = "C1"
$att.Cells.Item(699,5).Value = "221675@med.asu.edu.eg"
$att.Cells.Item(699,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(699,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(699,8).Value = "MICROBIOLOGY"
$att.Cells.Item(699,9).Value = "24/11/2025"
$att.Cells.Item(699,10).Value = "09:23:20"
$att.Cells.Item(699,11).Value = "C1"

# Row 700 - student 221527
$c = $att.Cells.Item(700,1)
$c.NumberFormat = "@"
$c.Value = "221527"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(700,2).Value = "شهد محمد عبدالرحمن ادريس"
$att.Cells.Item(700,3).Value = "Year 2"
$att.Cells.Item(700,4).Value = "C1"
$att.Cells.Item(700,5).Value = "221527@med.asu.edu.eg"
$att.Cells.Item(700,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(700,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(700,8).Value = "MICROBIOLOGY"
$att.Cells.Item(700,9).Value = "24/11/2025"
$att.Cells.Item(700,10).Value = "09:23:44"
$att.Cells.Item(700,11).Value = "C1"

# Row 701 - student 211704
$c = $att.Cells.Item(701,1)
$c.NumberFormat = "@"
$c.Value = "211704"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(701,2).Value = "الياس احمد بكردان"
$att.Cells.Item(701,3).Value = "Year 2"
$att.Cells.Item(701,4).Value = "C1"
$att.Cells.Item(701,5).Value = "211704@med.asu.edu.eg"
$att.Cells.Item(701,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(701,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(701,8).Value = "MICROBIOLOGY"
$att.Cells.Item(701,9).Value = "24/11/2025"
$att.Cells.Item(701,10).Value = "09:23:51"
$att.Cells.Item(701,11).Value = "C1"

# Row 702 - student 221324
$c = $att.Cells.Item(702,1)
$c.NumberFormat = "@"
$c.Value = "221324"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(702,2).Value = "مجد ذوقان خليل قيشاوي"
$att.Cells.Item(702,3).Value = "Year 2"
$att.Cells.Item(702,4).Value = "C1"
$att.Cells.Item(702,5).Value = "221324@med.asu.edu.eg"
$att.Cells.Item(702,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(702,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(702,8).Value = "MICROBIOLOGY"
$att.Cells.Item(702,9).Value = "24/11/2025"
$att.Cells.Item(702,10).Value = "09:23:57"
$att.Cells.Item(702,11).Value = "C1"

# Row 703 - student 221566
$c = $att.Cells.Item(703,1)
$c.NumberFormat = "@"
$c.Value = "221566"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(703,2).Value = "مصطفى سامى محمد عبد الله"
$att.Cells.Item(703,3).Value = "Year 2"
$att.Cells.Item(703,4).Value = "C1"
$att.Cells.Item(703,5).Value = "221566@med.asu.edu.eg"
$att.Cells.Item(703,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(703,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(703,8).Value = "MICROBIOLOGY"
$att.Cells.Item(703,9).Value = "24/11/2025"
$att.Cells.Item(703,10).Value = "09:24:09"
$att.Cells.Item(703,11).Value = "C1"

# Row 704 - student 221569
$c = $att.Cells.Item(704,1)
$c.NumberFormat = "@"
$c.Value = "221569"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(704,2).Value = "هبه جعفر محمد شوكت"
$att.Cells.Item(704,3).Value = "Year 2"
$att.Cells.Item(704,4).Value = "C1"
$att.Cells.Item(704,5).Value = "221569@med.asu.edu.eg"
$att.Cells.Item(704,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(704,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(704,8).Value = "MICROBIOLOGY"
$att.Cells.Item(704,9).Value = "24/11/2025"
$att.Cells.Item(704,10).Value = "09:24:13"
$att.Cells.Item(704,11).Value = "C1"

# Row 705 - student 212163
$c = $att.Cells.Item(705,1)
$c.NumberFormat = "@"
$c.Value = "212163"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(705,2).Value = "رقيه احمد عبد الله"
$att.Cells.Item(705,3).Value = "Year 2"
$att.Cells.Item(705,4).Value = "C1"
$att.Cells.Item(705,5).Value = "212163@med.asu.edu.eg"
$att.Cells.Item(705,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(705,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(705,8).Value = "MICROBIOLOGY"
$att.Cells.Item(705,9).Value = "24/11/2025"
$att.Cells.Item(705,10).Value = "09:24:52"
$att.Cells.Item(705,11).Value = "C1"

# Row 706 - student 221996
$c = $att.Cells.Item(706,1)
$c.NumberFormat = "@"
$c.Value = "221996"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(706,2).Value = "نياقوط فال توت دوير"
$att.Cells.Item(706,3).Value = "Year 2"
$att.Cells.Item(706,4).Value = "C1"
$att.Cells.Item(706,5).Value = "221996@med.asu.edu.eg"
$att.Cells.Item(706,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(706,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(706,8).Value = "MICROBIOLOGY"
$att.Cells.Item(706,9).Value = "24/11/2025"
$att.Cells.Item(706,10).Value = "09:25:08"
$att.Cells.Item(706,11).Value = "C1"

# Row 707 - student 221558
$c = $att.Cells.Item(707,1)
$c.NumberFormat = "@"
$c.Value = "221558"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(707,2).Value = "محمد عادل عوض باحاج"
$att.Cells.Item(707,3).Value = "Year 2"
$att.Cells.Item(707,4).Value = "C1"
$att.Cells.Item(707,5).Value = "221558@med.asu.edu.eg"
$att.Cells.Item(707,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(707,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(707,8).Value = "MICROBIOLOGY"
$att.Cells.Item(707,9).Value = "24/11/2025"
$att.Cells.Item(707,10).Value = "09:25:20"
$att.Cells.Item(707,11).Value = "C1"

# Row 708 - student 221494
$c = $att.Cells.Item(708,1)
$c.NumberFormat = "@"
$c.Value = "221494"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(708,2).Value = "حسن الصادق مصطفى الحاج"
$att.Cells.Item(708,3).Value = "Year 2"
$att.Cells.Item(708,4).Value = "C1"
$att.Cells.Item(708,5).Value = "221494@med.asu.edu.eg"
$att.Cells.Item(708,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(708,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(708,8).Value = "MICROBIOLOGY"
$att.Cells.Item(708,9).Value = "24/11/2025"
$att.Cells.Item(708,10).Value = "09:25:24"
$att.Cells.Item(708,11).Value = "C1"

# Row 709 - student 221459
$c = $att.Cells.Item(709,1)
$c.NumberFormat = "@"
$c.Value = "221459"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(709,2).Value = "محمد الطيب محمد زين"
$att.Cells.Item(709,3).Value = "Year 2"
$att.Cells.Item(709,4).Value = "C1"
$att.Cells.Item(709,5).Value = "221459@med.asu.edu.eg"
$att.Cells.Item(709,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(709,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(709,8).Value = "MICROBIOLOGY"
$att.Cells.Item(709,9).Value = "24/11/2025"
$att.Cells.Item(709,10).Value = "09:25:36"
$att.Cells.Item(709,11).Value = "C1"

# Row 710 - student 221536
$c = $att.Cells.Item(710,1)
$c.NumberFormat = "@"
$c.Value = "221536"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(710,2).Value = "عبده دفع الله سليمان كوكو"
$att.Cells.Item(710,3).Value = "Year 2"
$att.Cells.Item(710,4).Value = "C1"
$att.Cells.Item(710,5).Value = "221536@med.asu.edu.eg"
$att.Cells.Item(710,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(710,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(710,8).Value = "MICROBIOLOGY"
$att.Cells.Item(710,9).Value = "24/11/2025"
$att.Cells.Item(710,10).Value = "09:25:44"
$att.Cells.Item(710,11).Value = "C1"

# Row 711 - student 221522
$c = $att.Cells.Item(711,1)
$c.NumberFormat = "@"
$c.Value = "221522"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(711,2).Value = "آدم محمد احمد البديرات"
$att.Cells.Item(711,3).Value = "Year 2"
$att.Cells.Item(711,4).Value = "C1"
$att.Cells.Item(711,5).Value = "221522@med.asu.edu.eg"
$att.Cells.Item(711,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(711,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(711,8).Value = "MICROBIOLOGY"
$att.Cells.Item(711,9).Value = "24/11/2025"
$att.Cells.Item(711,10).Value = "09:25:51"
$att.Cells.Item(711,11).Value = "C1"

# Row 712 - student 221546
$c = $att.Cells.Item(712,1)
$c.NumberFormat = "@"
$c.Value = "221546"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(712,2).Value = "محمدزين ابوبكر محمد زين احمد"
$att.Cells.Item(712,3).Value = "Year 2"
$att.Cells.Item(712,4).Value = "C1"
$att.Cells.Item(712,5).Value = "221546@med.asu.edu.eg"
$att.Cells.Item(712,6).Value = "MICROBIOLOGY"
$c = $att.Cells.Item(712,7)
$c.NumberFormat = "@"
$c.Value = "1"
$plainTextSrc.Copy()
$c.PasteSpecial(-4122)
$att.Cells.Item(712,8).Value = "MICROBIOLOGY"
$att.Cells.Item(712,9).Value = "24/11/2025"
$att.Cells.Item(712,10).Value = "09:26:09"
$att.Cells.Item(712,11).Value = "C1"

$excel.CutCopyMode = $false
